$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(5, 2, 1, 18),
    @(6, 8, 1, 12),
    @(5, 4, 3, 16),
    @(2, 8, 3, 12),
    @(4, 12, 6, 8),
    @(5, 16, 7, 4),
    @(5, 8, 4, 12),
    @(3, 6, 4, 14),
    @(4, 12, 3, 8),
    @(6, 12, 6, 8),
    @(2, 14, 3, 6),
    @(7, 15, 6, 5),
    @(3, 8, 4, 12),
    @(8, 18, 5, 2)
)

$startRow = 1184
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $vals = $data[$i]
    $ws.Cells.Item($row, 1).Value = $vals[0]
    $ws.Cells.Item($row, 2).Value = $vals[1]
    $ws.Cells.Item($row, 3).Value = $vals[2]
    $ws.Cells.Item($row, 4).Value = $vals[3]
}

$excel.ActiveWindow.ScrollRow = 1172
$ws.Range("M1192").Select()
